# Autogenerated on Wed Apr 01 2015 00:15:40 GMT+0000 (Coordinated Universal Time)
# Refresh the "Spain" MSME summary figures with more precise decimal values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: Enterprises density (per 1000 people) -> Micro / SMEs / MSMEs
$ws.Range("B13").Value = "65.06"
$ws.Range("C13").Value = "2.86"
$ws.Range("D13").Value = "67.92"

# Row 14: Employment (% of total) -> Micro / SMEs / MSMEs
$ws.Range("B14").Value = "30.22"
$ws.Range("C14").Value = "32.62"
$ws.Range("D14").Value = "62.84"

# Row 16: Enterprises (% of total) -> Micro / SMEs / MSMEs
$ws.Range("B16").Value = "93.79"
$ws.Range("C16").Value = "6.08"
$ws.Range("D16").Value = "99.88"
